# Auto-generated Excel COM-interop script to apply market price / profit updates
# to the Kujata_Profits workbook sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1557.6666
$ws.Range("I28").Value = 1830.75
$ws.Range("J28").Value = 192.25
$ws.Range("K28").Value = 1830.75
$ws.Range("L28").Value = 192.25
$ws.Range("M28").Value = -1345.75
$ws.Range("N28").Value = -1162.25

$ws.Range("H33").Value = 261.3793
$ws.Range("I33").Value = 195.08696
$ws.Range("J33").Value = 515.5
$ws.Range("K33").Value = 195.08696
$ws.Range("L33").Value = 515.5
$ws.Range("M33").Value = 33.91304
$ws.Range("N33").Value = -973.5

$ws.Range("H43").Value = 3094070.2
$ws.Range("I43").Value = 11082.091
$ws.Range("J43").Value = 7938766
$ws.Range("K43").Value = 11082.091
$ws.Range("L43").Value = 7938766
$ws.Range("M43").Value = -11013.091
$ws.Range("N43").Value = -7938904

$ws.Range("H53").Value = 999.4400000000001
$ws.Range("I53").Value = 1232.2
$ws.Range("J53").Value = 68.40000000000001
$ws.Range("K53").Value = 1232.2
$ws.Range("L53").Value = 68.40000000000001
$ws.Range("M53").Value = -595.2
$ws.Range("N53").Value = -1342.4

$ws.Range("H64").Value = 5475
$ws.Range("J64").Value = 3960
$ws.Range("L64").Value = 3960
$ws.Range("N64").Value = -4456

$ws.Range("H67").Value = 5475
$ws.Range("J67").Value = 3960
$ws.Range("L67").Value = 3960
$ws.Range("N67").Value = -5676

$ws.Range("H69").Value = 3960
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 3960
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 11880
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -13628

$ws.Range("H72").Value = 3960
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3960
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 35640
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -44376

$ws.Range("H112").Value = 2892.1052
$ws.Range("J112").Value = 2991.6667
$ws.Range("L112").Value = 8975.000100000001
$ws.Range("N112").Value = -11191.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 55557984
$ws.Range("I61").Value = 58825924
$ws.Range("J61").Value = 2999
$ws.Range("K61").Value = 58825924
$ws.Range("L61").Value = 2999
$ws.Range("M61").Value = -58825712
$ws.Range("N61").Value = -3423

$ws.Range("H74").Value = 2233.4167
$ws.Range("I74").Value = 1780.1
$ws.Range("K74").Value = 1780.1
$ws.Range("M74").Value = -906.0999999999999

$ws.Range("H77").Value = 2233.4167
$ws.Range("I77").Value = 1780.1
$ws.Range("K77").Value = 8900.5
$ws.Range("M77").Value = -4532.5

$ws.Range("H132").Value = 2479.6538
$ws.Range("I132").Value = 1970.8422
$ws.Range("K132").Value = 5912.5266
$ws.Range("M132").Value = -3382.5266

$ws.Range("H136").Value = 55557984
$ws.Range("I136").Value = 58825924
$ws.Range("J136").Value = 2999
$ws.Range("K136").Value = 176477772
$ws.Range("L136").Value = 8997
$ws.Range("M136").Value = -176475222
$ws.Range("N136").Value = -14097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6359.091
$ws.Range("I134").Value = 1146.8462
$ws.Range("J134").Value = 13887.889
$ws.Range("K134").Value = 3440.5386
$ws.Range("L134").Value = 41663.667
$ws.Range("M134").Value = -905.5385999999999
$ws.Range("N134").Value = -46733.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 362.22223
$ws.Range("I22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = 150

$ws.Range("H31").Value = 1266.5476
$ws.Range("I31").Value = 1096.8286
$ws.Range("J31").Value = 2115.1428
$ws.Range("K31").Value = 1096.8286
$ws.Range("L31").Value = 2115.1428
$ws.Range("M31").Value = -801.8286000000001
$ws.Range("N31").Value = -2705.1428

$ws.Range("H34").Value = 1266.5476
$ws.Range("I34").Value = 1096.8286
$ws.Range("J34").Value = 2115.1428
$ws.Range("K34").Value = 1096.8286
$ws.Range("L34").Value = 2115.1428
$ws.Range("M34").Value = -894.8286000000001
$ws.Range("N34").Value = -2519.1428

$ws.Range("H58").Value = 1414.4062
$ws.Range("I58").Value = 991.2857
$ws.Range("J58").Value = 2222.182
$ws.Range("K58").Value = 991.2857
$ws.Range("L58").Value = 2222.182
$ws.Range("M58").Value = -788.2857
$ws.Range("N58").Value = -2628.182

$ws.Range("H107").Value = 2200
$ws.Range("J107").Value = 2733.3333
$ws.Range("L107").Value = 2733.3333
$ws.Range("N107").Value = -6573.3333

$ws.Range("H132").Value = 19637
$ws.Range("I132").Value = 34941.332
$ws.Range("K132").Value = 104823.996
$ws.Range("M132").Value = -102293.996

$ws.Range("H134").Value = 25002626
$ws.Range("I134").Value = 2907.8125
$ws.Range("K134").Value = 8723.4375
$ws.Range("M134").Value = -6188.4375

$ws.Range("H136").Value = 1414.4062
$ws.Range("I136").Value = 991.2857
$ws.Range("J136").Value = 2222.182
$ws.Range("K136").Value = 2973.8571
$ws.Range("L136").Value = 6666.545999999999
$ws.Range("M136").Value = -423.8571000000002
$ws.Range("N136").Value = -11766.546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 5000
$ws.Range("J23").Value = 5000
$ws.Range("L23").Value = 5000
$ws.Range("N23").Value = -5446

$ws.Range("H29").Value = 7800
$ws.Range("J29").Value = 7800
$ws.Range("L29").Value = 7800
$ws.Range("N29").Value = -8380

$ws.Range("H80").Value = 3518.0908
$ws.Range("I80").Value = 2649.5
$ws.Range("J80").Value = 3711.111
$ws.Range("K80").Value = 2649.5
$ws.Range("L80").Value = 3711.111
$ws.Range("M80").Value = -1651.5
$ws.Range("N80").Value = -5707.111

$ws.Range("H83").Value = 3518.0908
$ws.Range("I83").Value = 2649.5
$ws.Range("J83").Value = 3711.111
$ws.Range("K83").Value = 13247.5
$ws.Range("L83").Value = 18555.555
$ws.Range("M83").Value = -8255.5
$ws.Range("N83").Value = -28539.555

$ws.Range("H86").Value = 22497
$ws.Range("J86").Value = 22497
$ws.Range("L86").Value = 22497
$ws.Range("N86").Value = -24869

$ws.Range("H89").Value = 22497
$ws.Range("J89").Value = 22497
$ws.Range("L89").Value = 67491
$ws.Range("N89").Value = -79347

$ws.Range("H113").Value = 1030.375
$ws.Range("I113").Value = 1012.8889
$ws.Range("J113").Value = 1052.8572
$ws.Range("K113").Value = 1012.8889
$ws.Range("L113").Value = 1052.8572
$ws.Range("M113").Value = 1157.1111
$ws.Range("N113").Value = -5392.8572

$ws.Range("H122").Value = 1827.8462
$ws.Range("I122").Value = 1568
$ws.Range("J122").Value = 2412.5
$ws.Range("K122").Value = 4704
$ws.Range("L122").Value = 7237.5
$ws.Range("M122").Value = -2254
$ws.Range("N122").Value = -12137.5

$ws.Range("H132").Value = 4818.091
$ws.Range("I132").Value = 5000.125
$ws.Range("K132").Value = 15000.375
$ws.Range("M132").Value = -12470.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2954.2173
$ws.Range("I40").Value = 1969.6666
$ws.Range("K40").Value = 1969.6666
$ws.Range("M40").Value = -1833.6666

$ws.Range("H132").Value = 82320.64
$ws.Range("I132").Value = 17212.285
$ws.Range("J132").Value = 147429
$ws.Range("K132").Value = 51636.855
$ws.Range("L132").Value = 442287
$ws.Range("M132").Value = -49106.855
$ws.Range("N132").Value = -447347

$ws.Range("H136").Value = 20360.6
$ws.Range("I136").Value = 20360.6
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 61081.8
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -58531.8
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 60019
$ws.Range("J28").Value = 60019
$ws.Range("L28").Value = 60019
$ws.Range("N28").Value = -60715

$ws.Range("H53").Value = 12000
$ws.Range("J53").Value = 12000
$ws.Range("L53").Value = 12000
$ws.Range("N53").Value = -13214

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H101").Value = 9449.5
$ws.Range("J101").Value = 9449.5
$ws.Range("L101").Value = 9449.5
$ws.Range("N101").Value = -15939.5

$ws.Range("H122").Value = 10003715
$ws.Range("I122").Value = 10420495
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 31261485
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -31259035
$ws.Range("N122").Value = -7900
